$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The cell E67 (recovered count for 18 May 2020) was corrected from the
# number 1825 to the text "o" in the source data.
$ws.Range("E67").Value = "o"

# Add the new row of data for 19 Mayis 2020 (19 May 2020).
$ws.Range("A69").Value = 43970
$ws.Range("B69").Value = 25382
$ws.Range("C69").Value = 1022
$ws.Range("D69").Value = 28
$ws.Range("E69").Value = 1318

# Grow the worksheet table (Table3) so it covers the newly added row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E69"))

# Move the active selection to the new last populated row in column E,
# matching where the user's cursor ended up after the edit.
$ws.Range("E68").Select()
